$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.472.18'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.701.39'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.78'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5482'
$ws.Range('E6').Value = '  +4.40%  '
$ws.Range('E7').Value = '  +0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2752'
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('E10').Value = '  +0.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07703'
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('D12').Value = '1.695.87'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.563'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008397'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.99'
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').Value = '26.525.33'
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.957'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.01'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.83'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.269'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.30'
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1326'
$ws.Range('E25').Value = '  +7.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.925'
$ws.Range('E26').Value = '  +2.87%  '
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06252'
$ws.Range('E28').Value = '  -6.06%  '
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.621'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.609'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.043'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6190'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.415'
$ws.Range('E36').Value = '  +0.74%  '
$ws.Range('E37').Value = '  +2.52%  '
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('D39').Value = '1.121.23'
$ws.Range('E39').Value = '  +1.49%  '
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8801'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.017'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.47'
$ws.Range('D44').Value = '1.854.94'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.70'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.264'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.150'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('E51').Value = '  +0.10%  '
